$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 442 (shifts existing rows 442:469 down to 443:470)
$ws.Rows.Item(442).Insert()

# Populate the newly inserted row 442 with the new record
$ws.Range("A442").Value = 3
$ws.Range("B442").Value = "Femacal de La Calera"
$ws.Range("C442").Value = "Coquimbo"
$ws.Range("D442").Value = 44706
$ws.Range("E442").Value = 5
$ws.Range("F442").Value = 100112032
$ws.Range("G442").Value = "Zapallo italiano"
$ws.Range("H442").Value = "Sin especificar"
$ws.Range("I442").Value = "Primera"
$ws.Range("J442").Value = 215
$ws.Range("K442").Value = 12000
$ws.Range("L442").Value = 13000
$ws.Range("M442").Value = 12488
$ws.Range("N442").Value = "$/caja 70 unidades"
$ws.Range("O442").Value = "Región de Arica y Parinacota"
$ws.Range("P442").Value = 178
$ws.Range("Q442").Value = 70
$ws.Range("R442").Value = "Hortaliza"
